# Metagenomics - Re-analysis: add worked example values to the
# "Metagenomics" data sheet (2nd worksheet), row 2 of the annotation table.
#
# Commit message: "add example values to metagenomics templates"
#
# The example is the mangrove-sediment metagenome re-analysis study:
#   Andreote et al., "The Microbiome of Brazilian Mangrove Sediments as
#   Revealed by Metagenomics", PLoS ONE (2012), doi:10.1371/journal.pone.0038600,
#   MG-RAST id mgm4451033.3, sampled 2011-08-20 near Bertioga/SP, Brazil.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(2)

# --- Publication / accession columns -------------------------------------
$ws.Range("B2").Value = "The Microbiome of Brazilian Mangrove Sediments as Revealed by Metagenomics"
$ws.Range("E2").Value = "10.1371/journal.pone.0038600"
$ws.Range("K2").Value = "mgm4451033.3"

# --- Organism ---------------------------------------------------------------
$ws.Range("N2").Value = "mangrove metagenome"

# --- Term Source REF / Term Accession Number for Organism (stored as text) -
$ws.Range("O2:Q2").NumberFormat = "@"
$ws.Range("O2").Value = "NCBITaxon"
$ws.Range("P2").Value = "http://purl.obolibrary.org/obo/NCBITaxon_1284368"
$ws.Range("Q2").Value = "2011-08-20"

# --- biome -------------------------------------------------------------------
$ws.Range("T2").Value = "mangrove biome"

# --- environment (feature) / environmental material / geographic location /
#     GPS coordinates blocks (also stored as text) --------------------------
$ws.Range("U2:AE2").NumberFormat = "@"
$ws.Range("U2").Value = "ENVO"
$ws.Range("V2").Value = "http://purl.obolibrary.org/obo/ENVO_01000181"
$ws.Range("W2").Value = "tropical mangrove"
$ws.Range("X2").Value = "ENVO"
$ws.Range("Y2").Value = "http://purl.obolibrary.org/obo/ENVO_01000403"
$ws.Range("Z2").Value = "estuarine mud"
$ws.Range("AA2").Value = "ENVO"
$ws.Range("AB2").Value = "http://purl.obolibrary.org/obo/ENVO_00002160"
$ws.Range("AC2").Value = "Brazil"
$ws.Range("AD2").Value = "NCIT"
$ws.Range("AE2").Value = "http://purl.obolibrary.org/obo/NCIT_C16364"

$ws.Range("AF2").Value = "14° 14´6.0144""S, 51°55´31.0152""W"

# --- Depth parameter is numeric (0.3 meter) ---------------------------------
$ws.Range("AI2").Value = 0.3

# --- Column widths for the two columns whose example text is now the
#     longest entry in the table (Organism term source / Collection Date) --
$ws.Columns.Item(14).EntireColumn.AutoFit()
$ws.Columns.Item(17).EntireColumn.AutoFit()

# --- Selection / active sheet state -----------------------------------------
$ws.Activate()
$ws.Range("AJ7").Select()
